$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 244.5
$ws.Range("I2").Value = 250.57143
$ws.Range("J2").Value = 202
$ws.Range("K2").Value = 250.57143
$ws.Range("L2").Value = 202
$ws.Range("M2").Value = -137.57143
$ws.Range("N2").Value = -428
$ws.Range("H12").Value = 492.15384
$ws.Range("I12").Value = 509
$ws.Range("J12").Value = 290
$ws.Range("K12").Value = 509
$ws.Range("L12").Value = 290
$ws.Range("M12").Value = -339
$ws.Range("N12").Value = -630
$ws.Range("H51").Value = 7419.8
$ws.Range("I51").Value = 6474.75
$ws.Range("J51").Value = 8049.8335
$ws.Range("K51").Value = 6474.75
$ws.Range("L51").Value = 8049.8335
$ws.Range("M51").Value = -5990.75
$ws.Range("N51").Value = -9017.833500000001
$ws.Range("H132").Value = 13273.6
$ws.Range("I132").Value = 15342
$ws.Range("K132").Value = 46026
$ws.Range("M132").Value = -43496

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 841
$ws.Range("I4").Value = 799.5
$ws.Range("K4").Value = 799.5
$ws.Range("M4").Value = -683.5
$ws.Range("H17").Value = 9299
$ws.Range("I17").Value = 3899
$ws.Range("K17").Value = 3899
$ws.Range("M17").Value = -3726
$ws.Range("H32").Value = 3537.476
$ws.Range("I32").Value = 2125.6316
$ws.Range("K32").Value = 2125.6316
$ws.Range("M32").Value = -1838.6316
$ws.Range("H36").Value = 2821.5
$ws.Range("I36").Value = 2821.5
$ws.Range("K36").Value = 2821.5
$ws.Range("M36").Value = -2475.5
$ws.Range("H61").Value = 1636.3636
$ws.Range("J61").Value = 3000
$ws.Range("L61").Value = 3000
$ws.Range("N61").Value = -3424
$ws.Range("H122").Value = 1716.6897
$ws.Range("I122").Value = 1825.7693
$ws.Range("J122").Value = 771.3333
$ws.Range("K122").Value = 5477.3079
$ws.Range("L122").Value = 2313.9999
$ws.Range("M122").Value = -3027.3079
$ws.Range("N122").Value = -7213.9999
$ws.Range("H136").Value = 1636.3636
$ws.Range("J136").Value = 3000
$ws.Range("L136").Value = 9000
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2147.6667
$ws.Range("I20").Value = 2147.6667
$ws.Range("K20").Value = 2147.6667
$ws.Range("M20").Value = -1900.6667
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H80").Value = 2223.4443
$ws.Range("I80").Value = 373.6
$ws.Range("K80").Value = 373.6
$ws.Range("M80").Value = 624.4
$ws.Range("H83").Value = 2223.4443
$ws.Range("I83").Value = 373.6
$ws.Range("K83").Value = 1868
$ws.Range("M83").Value = 3124
$ws.Range("H86").Value = 2965.8
$ws.Range("I86").Value = 3184.2222
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 3184.2222
$ws.Range("L86").Value = 1000
$ws.Range("M86").Value = -2061.2222
$ws.Range("N86").Value = -3246
$ws.Range("H89").Value = 2965.8
$ws.Range("I89").Value = 3184.2222
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 15921.111
$ws.Range("L89").Value = 5000
$ws.Range("M89").Value = -10305.111
$ws.Range("N89").Value = -16232

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 7833.3335
$ws.Range("J9").Value = 8264.706
$ws.Range("L9").Value = 24794.118
$ws.Range("N9").Value = -25242.118
$ws.Range("H12").Value = 1256.6428
$ws.Range("J12").Value = 2020.375
$ws.Range("L12").Value = 6061.125
$ws.Range("N12").Value = -6407.125
$ws.Range("H107").Value = 558.6923
$ws.Range("I107").Value = 745.75
$ws.Range("J107").Value = 475.55554
$ws.Range("K107").Value = 2237.25
$ws.Range("L107").Value = 1426.66662
$ws.Range("M107").Value = -317.25
$ws.Range("N107").Value = -5266.66662
$ws.Range("H137").Value = 4903.143
$ws.Range("I137").Value = 3628.3333
$ws.Range("J137").Value = 5250.8184
$ws.Range("K137").Value = 10884.9999
$ws.Range("L137").Value = 15752.4552
$ws.Range("M137").Value = -5784.999899999999
$ws.Range("N137").Value = -25952.4552

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 6177.6
$ws.Range("I5").Value = 6177.6
$ws.Range("K5").Value = 6177.6
$ws.Range("M5").Value = -6065.6
$ws.Range("H26").Value = 41519.5
$ws.Range("J26").Value = 41519.5
$ws.Range("L26").Value = 41519.5
$ws.Range("N26").Value = -42079.5
$ws.Range("H33").Value = 25000
$ws.Range("J33").Value = 25000
$ws.Range("L33").Value = 25000
$ws.Range("N33").Value = -25504
$ws.Range("H50").Value = 41519.5
$ws.Range("J50").Value = 41519.5
$ws.Range("L50").Value = 41519.5
$ws.Range("N50").Value = -42515.5
$ws.Range("H55").Value = 13079.571
$ws.Range("I55").Value = 10686.667
$ws.Range("J55").Value = 14874.25
$ws.Range("K55").Value = 10686.667
$ws.Range("L55").Value = 14874.25
$ws.Range("M55").Value = -10359.667
$ws.Range("N55").Value = -15528.25
$ws.Range("H70").Value = 7839.909
$ws.Range("I70").Value = 6685
$ws.Range("J70").Value = 8499.857
$ws.Range("K70").Value = 6685
$ws.Range("L70").Value = 8499.857
$ws.Range("M70").Value = -6415
$ws.Range("N70").Value = -9039.857
$ws.Range("H73").Value = 7839.909
$ws.Range("I73").Value = 6685
$ws.Range("J73").Value = 8499.857
$ws.Range("K73").Value = 6685
$ws.Range("L73").Value = 8499.857
$ws.Range("M73").Value = -5749
$ws.Range("N73").Value = -10371.857
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -2830
$ws.Range("H122").Value = 2264
$ws.Range("I122").Value = 1696.7391
$ws.Range("J122").Value = 4127.857
$ws.Range("K122").Value = 5090.2173
$ws.Range("L122").Value = 12383.571
$ws.Range("M122").Value = -2640.2173
$ws.Range("N122").Value = -17283.571
$ws.Range("H126").Value = 2944.5
$ws.Range("J126").Value = 2944.5
$ws.Range("L126").Value = 8833.5
$ws.Range("N126").Value = -13773.5
$ws.Range("H132").Value = 1499.5
$ws.Range("I132").Value = 1499.5
$ws.Range("K132").Value = 4498.5
$ws.Range("M132").Value = -1968.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2130.261
$ws.Range("I7").Value = 2217.7222
$ws.Range("J7").Value = 1815.4
$ws.Range("K7").Value = 2217.7222
$ws.Range("L7").Value = 1815.4
$ws.Range("M7").Value = -2105.7222
$ws.Range("N7").Value = -2039.4
$ws.Range("H16").Value = 1249
$ws.Range("I16").Value = 1249
$ws.Range("K16").Value = 1249
$ws.Range("M16").Value = -1079
$ws.Range("H22").Value = 1854.579
$ws.Range("I22").Value = 1560.4615
$ws.Range("J22").Value = 2491.8333
$ws.Range("K22").Value = 1560.4615
$ws.Range("L22").Value = 2491.8333
$ws.Range("M22").Value = -1265.4615
$ws.Range("N22").Value = -3081.8333
$ws.Range("H27").Value = 1854.579
$ws.Range("I27").Value = 1560.4615
$ws.Range("J27").Value = 2491.8333
$ws.Range("K27").Value = 1560.4615
$ws.Range("L27").Value = 2491.8333
$ws.Range("M27").Value = -1453.4615
$ws.Range("N27").Value = -2705.8333
$ws.Range("H29").Value = 3500
$ws.Range("I29").Value = 2000
$ws.Range("J29").Value = 5000
$ws.Range("K29").Value = 2000
$ws.Range("L29").Value = 5000
$ws.Range("M29").Value = -1705
$ws.Range("N29").Value = -5590
$ws.Range("H46").Value = 3455.7144
$ws.Range("I46").Value = 3530
$ws.Range("J46").Value = 3400
$ws.Range("K46").Value = 3530
$ws.Range("L46").Value = 3400
$ws.Range("M46").Value = -3342
$ws.Range("N46").Value = -3776
$ws.Range("H82").Value = 2897.7144
$ws.Range("I82").Value = 3428.3333
$ws.Range("J82").Value = 2499.75
$ws.Range("K82").Value = 3428.3333
$ws.Range("L82").Value = 2499.75
$ws.Range("M82").Value = -3067.3333
$ws.Range("N82").Value = -3221.75
$ws.Range("H85").Value = 2897.7144
$ws.Range("I85").Value = 3428.3333
$ws.Range("J85").Value = 2499.75
$ws.Range("K85").Value = 3428.3333
$ws.Range("L85").Value = 2499.75
$ws.Range("M85").Value = -2180.3333
$ws.Range("N85").Value = -4995.75
$ws.Range("H115").Value = 40000
$ws.Range("J115").Value = 40000
$ws.Range("L115").Value = 40000
$ws.Range("N115").Value = -42350
$ws.Range("H122").Value = 3589.4614
$ws.Range("I122").Value = 3423.9092
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 10271.7276
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -7821.7276
$ws.Range("N122").Value = -18400
$ws.Range("H126").Value = 2130.261
$ws.Range("I126").Value = 2217.7222
$ws.Range("J126").Value = 1815.4
$ws.Range("K126").Value = 6653.1666
$ws.Range("L126").Value = 5446.200000000001
$ws.Range("M126").Value = -4183.1666
$ws.Range("N126").Value = -10386.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 1000
$ws.Range("J5").Value = 1000
$ws.Range("L5").Value = 1000
$ws.Range("N5").Value = -1224
$ws.Range("H81").Value = 4999.375
$ws.Range("I81").Value = 4499.2
$ws.Range("J81").Value = 5833
$ws.Range("K81").Value = 8998.4
$ws.Range("L81").Value = 11666
$ws.Range("M81").Value = -7937.4
$ws.Range("N81").Value = -13788
$ws.Range("H84").Value = 4999.375
$ws.Range("I84").Value = 4499.2
$ws.Range("J84").Value = 5833
$ws.Range("K84").Value = 44992
$ws.Range("L84").Value = 58330
$ws.Range("M84").Value = -39688
$ws.Range("N84").Value = -68938
$ws.Range("H107").Value = 3087.5334
$ws.Range("I107").Value = 1861.3
$ws.Range("K107").Value = 5583.9
$ws.Range("M107").Value = -3663.9
$ws.Range("H136").Value = 2693.389
$ws.Range("I136").Value = 1969.3529
$ws.Range("K136").Value = 5908.0587
$ws.Range("M136").Value = -3358.0587
